$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.)
$newRows = @(
    @{ Row = 252; Date = 44326; B = 0; C = 8; D = 113.7980085348507 },
    @{ Row = 253; Date = 44327; B = 0; C = 8; D = 113.7980085348507 },
    @{ Row = 254; Date = 44328; B = 0; C = 7; D = 99.5732574679943 },
    @{ Row = 255; Date = 44329; B = 0; C = 4; D = 56.89900426742533 }
)

$lastRow = 251

foreach ($item in $newRows) {
    $r = $item.Row

    # Copy formatting (incl. date number format / borders / font) from the row above
    $ws.Range("A$lastRow`:D$lastRow").Copy()
    $ws.Range("A$r`:D$r").PasteSpecial(-4122)

    $ws.Range("A$r").Value = $item.Date
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D

    $lastRow = $r
}

$excel.CutCopyMode = 0
